Write-Output ($ppt | Get-Member | Out-String)
